# Rename the "tablets" sheet to "data"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "data"

# Replace the data block (A2:C11) with the new (rescaled) values
$data = @(
    @(51, 42, 72),
    @(78, 47, 88),
    @(28, 39, 44),
    @(39, 33, 46),
    @(50, 36, 65),
    @(34, 22, 45),
    @(26, 93, 37),
    @(25, 15, 36),
    @(19, 18, 22),
    @(44, 35, 53)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Give columns A:C an explicit (custom) width
$ws.Columns("A:C").ColumnWidth = 6.9

# Move the active selection to A13
$ws.Range("A13").Select()

# Configure page setup: A4 paper, portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
